$wb = $excel.ActiveWorkbook

# Rename the "Send_Default" sheet to "Test" (the remaining/consolidated Test sheet).
$wsTest = $wb.Worksheets.Item("Send_Default")
$wsTest.Name = "Test"

# Update the Config lookup table so both Send/Recv rows point at the single "Test" sheet.
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Range("D4").Value = "Test"
$wsConfig.Range("K4").Value = "Test"

# Make Config the active sheet/tab and set its selection.
$wsConfig.Activate()
$wsConfig.Range("K13").Select()

# Restore the (now inactive) Test sheet's own remembered selection.
$wsTest.Range("R26").Select()
